$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 24
$ws.Range("E2").Value = 22.53000068664551
$ws.Range("F2").Value = 29.44000053405762
$ws.Range("G2").Value = 18.89999961853028
$ws.Range("H2").Value = 1435682333
$ws.Range("I2").Value = "SNAP"

$ws.Range("D3").Value = 24
$ws.Range("E3").Value = 22.53000068664551
$ws.Range("F3").Value = 29.44000053405762
$ws.Range("G3").Value = 18.89999961853028
$ws.Range("H3").Value = 1435682333
$ws.Range("I3").Value = "SNAP"

$ws.Range("D4").Value = 24
$ws.Range("E4").Value = 22.53000068664551
$ws.Range("F4").Value = 29.44000053405762
$ws.Range("G4").Value = 18.89999961853028
$ws.Range("H4").Value = 1435682333
$ws.Range("I4").Value = "SNAP"

$ws.Range("D5").Value = 24
$ws.Range("E5").Value = 22.53000068664551
$ws.Range("F5").Value = 29.44000053405762
$ws.Range("G5").Value = 18.89999961853028
$ws.Range("H5").Value = 1435682333
$ws.Range("I5").Value = "SNAP"

$ws.Range("D6").Value = 24
$ws.Range("E6").Value = 22.53000068664551
$ws.Range("F6").Value = 29.44000053405762
$ws.Range("G6").Value = 18.89999961853028
$ws.Range("H6").Value = 1435682333
$ws.Range("I6").Value = "SNAP"

$ws.Range("D7").Value = 22.70000076293945
$ws.Range("E7").Value = 22.54999923706055
$ws.Range("F7").Value = 22.89999961853028
$ws.Range("G7").Value = 19.72999954223633
$ws.Range("H7").Value = 1435682333
$ws.Range("I7").Value = "SNAP"

$ws.Range("D8").Value = 17.90999984741211
$ws.Range("E8").Value = 13.67000007629394
$ws.Range("F8").Value = 17.92000007629395
$ws.Range("G8").Value = 13.10000038146973
$ws.Range("H8").Value = 1435682333
$ws.Range("I8").Value = "SNAP"

$ws.Range("D9").Value = 14.60000038146973
$ws.Range("E9").Value = 15.34000015258789
$ws.Range("F9").Value = 16.8799991607666
$ws.Range("G9").Value = 13.96000003814697
$ws.Range("H9").Value = 1435682333
$ws.Range("I9").Value = "SNAP"

$ws.Range("D10").Value = 14.6899995803833
$ws.Range("E10").Value = 13.52000045776367
$ws.Range("F10").Value = 15.38000011444092
$ws.Range("G10").Value = 13.14999961853027
$ws.Range("H10").Value = 1435682333
$ws.Range("I10").Value = "SNAP"

$ws.Range("D11").Value = 15.67000007629394
$ws.Range("E11").Value = 14.32999992370606
$ws.Range("F11").Value = 15.95699977874756
$ws.Range("G11").Value = 13.61999988555908
$ws.Range("H11").Value = 1435682333
$ws.Range("I11").Value = "SNAP"

$ws.Range("D12").Value = 13.05000019073486
$ws.Range("E12").Value = 12.5
$ws.Range("F12").Value = 14.13500022888184
$ws.Range("G12").Value = 12
$ws.Range("H12").Value = 1435682333
$ws.Range("I12").Value = "SNAP"

$ws.Range("D13").Value = 8.5
$ws.Range("E13").Value = 6.610000133514404
$ws.Range("F13").Value = 8.569999694824219
$ws.Range("G13").Value = 5.769999980926514
$ws.Range("H13").Value = 1435682333
$ws.Range("I13").Value = "SNAP"

$ws.Range("D14").Value = 5.380000114440918
$ws.Range("E14").Value = 6.679999828338623
$ws.Range("F14").Value = 6.809999942779541
$ws.Range("G14").Value = 5.349999904632568
$ws.Range("H14").Value = 1435682333
$ws.Range("I14").Value = "SNAP"

$ws.Range("D15").Value = 11.14000034332275
$ws.Range("E15").Value = 11.14000034332275
$ws.Range("F15").Value = 12.63000011444092
$ws.Range("G15").Value = 10.51000022888184
$ws.Range("H15").Value = 1435682333
$ws.Range("I15").Value = "SNAP"

$ws.Range("D17").Value = 16.04000091552734
$ws.Range("E17").Value = 15.0600004196167
$ws.Range("F17").Value = 16.45999908447266
$ws.Range("G17").Value = 12.71000003814697
$ws.Range("H17").Value = 1435682333
$ws.Range("I17").Value = "SNAP"

$ws.Range("D18").Value = 16.6200008392334
$ws.Range("E18").Value = 18.3799991607666
$ws.Range("F18").Value = 19.7549991607666
$ws.Range("G18").Value = 16.35000038146973
$ws.Range("H18").Value = 1435682333
$ws.Range("I18").Value = "SNAP"

$ws.Range("D19").Value = 11.31999969482422
$ws.Range("E19").Value = 17.61000061035156
$ws.Range("F19").Value = 17.98999977111816
$ws.Range("G19").Value = 10.85000038146973
$ws.Range("H19").Value = 1435682333
$ws.Range("I19").Value = "SNAP"

$ws.Range("D20").Value = 23.64999961853028
$ws.Range("E20").Value = 22.42000007629395
$ws.Range("F20").Value = 26.76000022888184
$ws.Range("G20").Value = 21.53000068664551
$ws.Range("H20").Value = 1435682333
$ws.Range("I20").Value = "SNAP"

$ws.Range("D21").Value = 26.28000068664551
$ws.Range("E21").Value = 39.38999938964844
$ws.Range("F21").Value = 44.18000030517578
$ws.Range("G21").Value = 26.1299991607666
$ws.Range("H21").Value = 1435682333
$ws.Range("I21").Value = "SNAP"

$ws.Range("D22").Value = 50.43999862670898
$ws.Range("E22").Value = 52.93999862670898
$ws.Range("F22").Value = 57.38999938964844
$ws.Range("G22").Value = 48.0989990234375
$ws.Range("H22").Value = 1435682333
$ws.Range("I22").Value = "SNAP"

$ws.Range("D23").Value = 53.68999862670898
$ws.Range("E23").Value = 61.81999969482422
$ws.Range("F23").Value = 65.86000061035156
$ws.Range("G23").Value = 52.68999862670898
$ws.Range("H23").Value = 1435682333
$ws.Range("I23").Value = "SNAP"

$ws.Range("D24").Value = 68.58000183105469
$ws.Range("E24").Value = 74.41999816894531
$ws.Range("F24").Value = 79.18000030517578
$ws.Range("G24").Value = 57.47999954223633
$ws.Range("H24").Value = 1435682333
$ws.Range("I24").Value = "SNAP"

$ws.Range("D25").Value = 74.81999969482422
$ws.Range("E25").Value = 52.58000183105469
$ws.Range("F25").Value = 79.30000305175781
$ws.Range("G25").Value = 51.65999984741211
$ws.Range("H25").Value = 1435682333
$ws.Range("I25").Value = "SNAP"

$ws.Range("D26").Value = 47.63000106811523
$ws.Range("E26").Value = 32.54000091552734
$ws.Range("F26").Value = 47.70999908447266
$ws.Range("G26").Value = 28.02000045776367
$ws.Range("H26").Value = 1435682333
$ws.Range("I26").Value = "SNAP"

$ws.Range("D27").Value = 36.20000076293945
$ws.Range("E27").Value = 28.45999908447266
$ws.Range("F27").Value = 39.79999923706055
$ws.Range("G27").Value = 26.44000053405762
$ws.Range("H27").Value = 1435682333
$ws.Range("I27").Value = "SNAP"

$ws.Range("D28").Value = 13.18000030517578
$ws.Range("E28").Value = 9.880000114440918
$ws.Range("F28").Value = 16.54500007629395
$ws.Range("G28").Value = 9.34000015258789
$ws.Range("H28").Value = 1435682333
$ws.Range("I28").Value = "SNAP"

$ws.Range("D29").Value = 9.960000038146973
$ws.Range("E29").Value = 9.90999984741211
$ws.Range("F29").Value = 11.85000038146973
$ws.Range("G29").Value = 7.329999923706055
$ws.Range("H29").Value = 1435682333
$ws.Range("I29").Value = "SNAP"

$ws.Range("D30").Value = 9.149999618530272
$ws.Range("E30").Value = 11.5600004196167
$ws.Range("F30").Value = 11.56999969482422
$ws.Range("G30").Value = 8.6899995803833
$ws.Range("H30").Value = 1435682333
$ws.Range("I30").Value = "SNAP"

$ws.Range("D31").Value = 11.06999969482422
$ws.Range("E31").Value = 8.710000038146973
$ws.Range("F31").Value = 11.47000026702881
$ws.Range("G31").Value = 8.40999984741211
$ws.Range("H31").Value = 1435682333
$ws.Range("I31").Value = "SNAP"

$ws.Range("D32").Value = 11.89999961853027
$ws.Range("E32").Value = 11.35999965667725
$ws.Range("F32").Value = 13.89000034332275
$ws.Range("G32").Value = 9.989999771118164
$ws.Range("H32").Value = 1435682333
$ws.Range("I32").Value = "SNAP"

$ws.Range("D33").Value = 8.869999885559082
$ws.Range("E33").Value = 10.01000022888184
$ws.Range("F33").Value = 10.86999988555908
$ws.Range("G33").Value = 8.395000457763672
$ws.Range("H33").Value = 1435682333
$ws.Range("I33").Value = "SNAP"

$ws.Range("D34").Value = 16.52000045776367
$ws.Range("E34").Value = 15.89000034332275
$ws.Range("F34").Value = 17.75
$ws.Range("G34").Value = 15.18000030517578
$ws.Range("H34").Value = 1435682333
$ws.Range("I34").Value = "SNAP"

$ws.Range("D35").Value = 11.47999954223633
$ws.Range("E35").Value = 15.05000019073486
$ws.Range("F35").Value = 15.35999965667725
$ws.Range("G35").Value = 10.08300018310547
$ws.Range("H35").Value = 1435682333
$ws.Range("I35").Value = "SNAP"

$ws.Range("D36").Value = 16.43000030517578
$ws.Range("E36").Value = 13.31999969482422
$ws.Range("F36").Value = 17.32999992370605
$ws.Range("G36").Value = 12.82999992370606
$ws.Range("H36").Value = 1435682333
$ws.Range("I36").Value = "SNAP"

$ws.Range("D37").Value = 10.77999973297119
$ws.Range("E37").Value = 12.15999984741211
$ws.Range("F37").Value = 12.82999992370606
$ws.Range("G37").Value = 9.9399995803833
$ws.Range("H37").Value = 1435682333
$ws.Range("I37").Value = "SNAP"

$ws.Range("D38").Value = 11
$ws.Range("E38").Value = 11.28999996185303
$ws.Range("F38").Value = 13.2810001373291
$ws.Range("G38").Value = 10.39999961853027
$ws.Range("H38").Value = 1435682333
$ws.Range("I38").Value = "SNAP"

$ws.Range("D39").Value = 8.739999771118164
$ws.Range("E39").Value = 7.960000038146973
$ws.Range("F39").Value = 9.229999542236328
$ws.Range("G39").Value = 7.079999923706055
$ws.Range("H39").Value = 1435682333
$ws.Range("I39").Value = "SNAP"

$ws.Range("D40").Value = 9
$ws.Range("E40").Value = 9.43000030517578
$ws.Range("F40").Value = 10.40999984741211
$ws.Range("G40").Value = 8.90999984741211
$ws.Range("H40").Value = 1435682333
$ws.Range("I40").Value = "SNAP"
